$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the rows we're about to edit are visible first, so Excel doesn't
# stamp an auto-fit row height on a hidden row when its cell values change
$ws.Rows(5).Hidden = $false
$ws.Rows(6).Hidden = $false
$ws.Rows(7).Hidden = $false

# Remove the stray note in A2 (its shared string is dropped entirely once unreferenced)
$ws.Range("A2").ClearContents()

# Rename the 4th series header from "KVM vhost" to "KVM"
$ws.Range("E4").Value() = "KVM"

# Row 5 (raw nuttcp numbers): update values and add the KVM column
$ws.Range("B5").Value() = 9348.1319
$ws.Range("C5").Value() = 9329.8451
$ws.Range("E5").Value() = 9294.2403

# Row 7 (raw numbers) gets new values and the KVM column
$ws.Range("B7").Value() = 9340.2296
$ws.Range("E7").Value() = 9264.9354

# Now apply the final row visibility: row 5 hides, row 6 (the /1000 view) shows,
# row 7 stays hidden, row 8 stays visible
$ws.Rows(5).Hidden = $true
$ws.Rows(7).Hidden = $true

# Move the active selection
$ws.Range("B13").Select() | Out-Null
